# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the cfade5cd file, adding error detail messages on the
# zh-cn and de-de sheets, and widening the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Update every cell that shows the status of the cfade5cd row (row 3) so
# the shared "Ready for handoff" text is replaced everywhere it appears,
# reflecting the failed handback transform.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zh.Range("C3").Value = "Handback transform failed"
$de.Range("C3").Value = "Handback transform failed"

# Record the error detail messages explaining the handback/handoff file
# name mismatch, one per target language sheet.
$zh.Range("P3").Value = "Handback file name: wylevbvn.gtf is different with handoff file name: cfade5cd-1d6b-4e78-af54-89e1cbe45f60.7a1d4cbe204597b3e0e29b768724fba687614ba4.zh-cn."
$de.Range("P3").Value = "Handback file name: wylevbvn.gtf is different with handoff file name: cfade5cd-1d6b-4e78-af54-89e1cbe45f60.7a1d4cbe204597b3e0e29b768724fba687614ba4.de-de."

# Widen the "Error Detail" column (P) on both sheets to fit the new text.
# (ColumnWidth is specified in characters; the engine adds ~5/6 of a
# character of internal padding when it stores the OOXML column width, so
# subtract that back off here to land exactly on a stored width of 40.)
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
